# Update cryptocurrency price/volume data in-place (canonical diff replay).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.427.41"
$ws.Range('E2').Value = '  -0.19%  '

$ws.Range('D3').Value = "'1.816.75"
$ws.Range('E3').Value = '  -0.65%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').Value = "'315.19"
$ws.Range('E5').Value = '  -0.70%  '

$ws.Range('E6').Value = '  +0.25%  '

$ws.Range('D7').Value = "'0.5084"
$ws.Range('E7').Value = '  -4.86%  '

$ws.Range('D8').Value = "'0.3954"
$ws.Range('E8').Value = '  -1.58%  '

$ws.Range('D9').Value = "'0.08224"
$ws.Range('E9').Value = '  +8.12%  '

$ws.Range('D10').Value = "'41.66"
$ws.Range('E10').Value = '  -0.44%  '

$ws.Range('D11').Value = "'1.106"
$ws.Range('E11').Value = '  -0.54%  '

$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('D13').Value = "'6.284"
$ws.Range('E13').Value = '  -0.90%  '

$ws.Range('E14').Value = '  +0.19%  '

$ws.Range('D15').Value = "'7.503"

$ws.Range('D16').Value = "'1.815.39"
$ws.Range('E16').Value = '  -0.70%  '

$ws.Range('D17').Value = "'0.00001140"
$ws.Range('E17').Value = '  +6.07%  '

$ws.Range('D18').Value = "'92.45"
$ws.Range('E18').Value = '  +3.23%  '

$ws.Range('E19').Value = '  +0.51%  '

$ws.Range('D20').Value = "'17.68"
$ws.Range('E20').Value = '  -0.20%  '

$ws.Range('E21').Value = '  +0.12%  '

$ws.Range('D22').Value = "'6.097"
$ws.Range('E22').Value = '  +0.07%  '

$ws.Range('D23').Value = "'28.454.33"
$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('D24').Value = "'11.28"
$ws.Range('E24').Value = '  +0.80%  '

$ws.Range('D25').Value = "'2.265"
$ws.Range('E25').Value = '  +2.34%  '

$ws.Range('D26').Value = "'21.17"
$ws.Range('E26').Value = '  +2.40%  '

$ws.Range('D27').Value = "'155.42"
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('D28').Value = "'2.026.29"
$ws.Range('E28').Value = '  -0.64%  '

$ws.Range('D29').Value = "'2.402"
$ws.Range('E29').Value = '  -2.43%  '

$ws.Range('D30').Value = "'126.03"
$ws.Range('E30').Value = '  +1.17%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'1.109"
$ws.Range('E31').Value = '  -1.33%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.1097"
$ws.Range('E32').Value = '  -0.55%  '

$ws.Range('D33').Value = "'5.780"
$ws.Range('E33').Value = '  +1.90%  '

$ws.Range('D34').Value = "'3.651"
$ws.Range('E34').Value = '  +0.28%  '

$ws.Range('D35').Value = "'0.07026"
$ws.Range('E35').Value = '  -7.08%  '

$ws.Range('D36').Value = "'0.2219"
$ws.Range('E36').Value = '  -0.51%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.02328"
$ws.Range('E37').Value = '  -0.66%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'5.212"
$ws.Range('E38').Value = '  +0.03%  '

$ws.Range('E39').Value = '  -0.40%  '

$ws.Range('D40').Value = "'0.6275"
$ws.Range('E40').Value = '  +0.11%  '

$ws.Range('D41').Value = "'11.26"
$ws.Range('E41').Value = '  -0.42%  '

$ws.Range('D42').Value = "'1.175"
$ws.Range('E42').Value = '  +0.05%  '

$ws.Range('E43').Value = '  +0.14%  '

$ws.Range('D44').Value = "'1.404"
$ws.Range('E44').Value = '  +0.88%  '

$ws.Range('D45').Value = "'13.41"
$ws.Range('E45').Value = '  -0.58%  '

$ws.Range('D46').Value = "'3.742"
$ws.Range('E46').Value = '  +1.05%  '

$ws.Range('D47').Value = "'0.5903"
$ws.Range('E47').Value = '  +0.90%  '

$ws.Range('D48').Value = "'124.95"
$ws.Range('E48').Value = '  +0.07%  '

$ws.Range('D49').Value = "'1.975"
$ws.Range('E49').Value = '  -1.15%  '

$ws.Range('D50').Value = "'1.185"
$ws.Range('E50').Value = '  -1.32%  '

$ws.Range('D51').Value = "'0.06888"
$ws.Range('E51').Value = '  +0.00%  '

